$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold plain text in this sheet
# (e.g. "1.002", "29.784.86", "  -0.38%  "), never real numbers/percents.
# Force text format on the cells we are about to (re)write so Excel's
# auto-detection doesn't silently convert decimal-looking strings
# ("1.003", "0.8800", "5.980", ...) into floating point numbers and mangle
# their textual representation (e.g. dropping trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

# Swap rows 13/14 (Polkadot <-> WrappedEther) and update values
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.865.61"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.308"
$ws.Range("E14").Value = "  -1.91%  "

# Swap rows 21/22 (WrappedliquidstakedEther2.0 <-> Dai) and update values
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.136.97"
$ws.Range("E22").Value = "  +0.50%  "

# Swap rows 46/47 (Quant <-> Maker) and update values
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.036.93"
$ws.Range("E46").Value = "  -6.79%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "103.48"
$ws.Range("E47").Value = "  -1.19%  "

# Simple value updates (D/E columns only) for all other rows
$updates = @{
    2  = @{ D = "29.764.71"; E = "  -0.37%  " }
    3  = @{ D = "1.871.55";  E = "  -0.07%  " }
    4  = @{ D = "1.003";     E = "  +0.49%  " }
    5  = @{ D = "0.7161";    E = "  -3.01%  " }
    6  = @{ D = "241.55";    E = "  -0.16%  " }
    7  = @{ D = "1.003";     E = "  +0.39%  " }
    8  = @{ D = "0.3143";    E = "  -0.36%  " }
    9  = @{ D = "0.07506";   E = "  +4.45%  " }
    10 = @{ D = "24.46";     E = "  -1.05%  " }
    11 = @{ D = "0.08181";   E = "  -2.09%  " }
    12 = @{ D = "0.7416";    E = "  -1.15%  " }
    15 = @{ D = "92.29";     E = "  -0.31%  " }
    16 = @{ D = "29.840.98"; E = "  -0.12%  " }
    17 = @{ D = "5.994";     E = "  -1.22%  " }
    18 = @{ D = "245.92";    E = "  -0.11%  " }
    19 = @{ D = "0.000007904"; E = "  +0.95%  " }
    20 = @{ D = "13.43";     E = "  -0.97%  " }
    23 = @{ D = "1.003";     E = "  +0.55%  " }
    24 = @{ D = "7.705";     E = "  -3.62%  " }
    25 = @{ D = "9.168";     E = "  -0.96%  " }
    26 = @{ D = "0.1492";    E = "  -3.87%  " }
    27 = @{ D = "163.55" }
    28 = @{ D = "18.52";     E = "  -0.68%  " }
    29 = @{ D = "1.998";     E = "  -1.62%  " }
    30 = @{ D = "1.423";     E = "  -5.46%  " }
    31 = @{ D = "4.523";     E = "  -1.28%  " }
    32 = @{ D = "1.524";     E = "  -0.72%  " }
    33 = @{ D = "4.168";     E = "  -2.41%  " }
    34 = @{ D = "0.05434";   E = "  +2.36%  " }
    35 = @{ D = "1.221";     E = "  -1.29%  " }
    36 = @{ D = "0.7356";    E = "  -2.47%  " }
    37 = @{ D = "0.9978";    E = "  -0.09%  " }
    38 = @{ D = "2.701";     E = "  +0.21%  " }
    39 = @{ D = "0.01907";   E = "  -2.63%  " }
    40 = @{ D = "2.729";     E = "  -0.82%  " }
    41 = @{ D = "0.4442";    E = "  -1.50%  " }
    42 = @{ D = "0.8800";    E = "  +3.03%  " }
    43 = @{ D = "5.980";     E = "  -1.07%  " }
    44 = @{ D = "71.18";     E = "  -1.50%  " }
    45 = @{ D = "1.002";     E = "  +0.32%  " }
    48 = @{ D = "7.442";     E = "  -2.41%  " }
    49 = @{ D = "9.557";     E = "  +0.72%  " }
    50 = @{ D = "1.785";     E = "  -2.97%  " }
    51 = @{ D = "2.033.25";  E = "  +0.53%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
